$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells keep their original text (non-numeric) formatting
# by forcing Text number format before assigning the new values, matching
# the source file's inline-string representation.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "302.63"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "5.24%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "34.82"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "12.33%"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "4.19%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07753"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "5.06%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.344"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "5.39%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.018"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "4.28%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.949"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "5.27%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9276"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.84%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1012"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "15.88%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1794"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "6.65%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08550"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "4.10%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03317"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "6.68%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09887"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.67%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001505"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.66%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005789"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.39%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.468"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.68%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.168"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "3.87%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.15%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1310"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.08%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.336"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "13.07%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2388"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "9.20%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04566"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.21%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001218"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.54%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004460"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "7.64%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001300"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.26%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-0.14%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01779"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "12.55%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04754"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "6.43%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007749"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "5.14%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "6.73%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007095"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-25.94%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002151"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "0.23%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009201"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "10.13%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006122"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.27%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.06%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.646"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "20.55%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002001"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.06%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002101"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.06%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002001"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.06%"
